$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Octubre de 2020 a las 02:37"

# 2. Swap the "Islas Malvinas" / "Montserrat" rows (they were re-ordered in the
#    shared-strings table, which - combined with the underlying data refresh -
#    results in the two countries exchanging their Casos activos (D) / Muertes (H)
#    values for rows 215 and 216).
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

$ws.Range("A216").Value = "Montserrat"
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1

# 3. Update the refreshed covid statistics for the affected countries.

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7635691
$ws.Range("C4").Value = 33226
$ws.Range("D4").Value = 4844656
$ws.Range("E4").Value = 2576426
$ws.Range("G4").Value = 330
$ws.Range("H4").Value = 214609

# Row 26 - Alemania
$ws.Range("B26").Value = 301571
$ws.Range("C26").Value = 1543
$ws.Range("E26").Value = 30069

# Row 50 - Chequia
$ws.Range("B50").Value = 82446
$ws.Range("C50").Value = 1841
$ws.Range("D50").Value = 44149
$ws.Range("E50").Value = 37570
$ws.Range("G50").Value = 16
$ws.Range("H50").Value = 727

# Row 69 - Paraguay
$ws.Range("B69").Value = 44182
$ws.Range("C69").Value = 730
$ws.Range("D69").Value = 27203
$ws.Range("E69").Value = 16050
$ws.Range("G69").Value = 16
$ws.Range("H69").Value = 929

# Row 130 - Surinam
$ws.Range("B130").Value = 4941
$ws.Range("C130").Value = 17
$ws.Range("D130").Value = 4741
$ws.Range("E130").Value = 94

# Row 147 - Guyana
$ws.Range("B147").Value = 3093
$ws.Range("C147").Value = 88
$ws.Range("D147").Value = 1920
$ws.Range("E147").Value = 1086
$ws.Range("G147").Value = 2
$ws.Range("H147").Value = 87

# Row 153 - Burkina Faso
$ws.Range("B153").Value = 2167
$ws.Range("C153").Value = 13
$ws.Range("D153").Value = 1419
$ws.Range("E153").Value = 689

# Row 172 - Islas Turcas y Caicos
$ws.Range("D172").Value = 656
$ws.Range("E172").Value = 33

# Row 191 - Barbados
$ws.Range("B191").Value = 199
$ws.Range("C191").Value = 3
$ws.Range("E191").Value = 10

$wb.Save()
